$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "[Lucio%Verdoni%NULL%1,                Angelo%Mazza%NULL%1,                Annalisa%Gervasoni%NULL%1,                Laura%Martelli%NULL%1,                Maurizio%Ruggeri%NULL%1,                Matteo%Ciuffreda%NULL%1,                Ezio%Bonanomi%NULL%1,                Lorenzo%D'Antiga%NULL%1]"
$ws.Range("E3").Value = "[Eva W.%Cheung%xref no email%1,        Philip%Zachariah%xref no email%1,        Mark%Gorelik%xref no email%1,        Alexis%Boneparth%xref no email%1,        Steven G.%Kernie%xref no email%1,        Jordan S.%Orange%xref no email%1,        Joshua D.%Milner%xref no email%1]"
$ws.Range("E5").Value = "[Shubhi%Kaushik%NULL%1,                Scott I.%Aydin%NULL%1,                Kim R.%Derespina%NULL%0,                Prerna B.%Bansal%NULL%1,                Shanna%Kowalsky%NULL%2,                Rebecca%Trachtman%NULL%2,                Jennifer K.%Gillen%NULL%1,                Michelle M.%Perez%NULL%1,                Sara H.%Soshnick%NULL%1,                Edward E.%Conway%NULL%1,                Asher%Bercow%NULL%1,                Howard S.%Seiden%NULL%1,                Robert H.%Pass%NULL%1,                Henry M.%Ushay%NULL%0,                George%Ofori-Amanfo%NULL%1,                Shivanand S.%Medar%NULL%0]"
$ws.Range("E6").Value = "[Jonathan%Miller%NULL%1,                Amanda%Cantor%NULL%1,                Philip%Zachariah%NULL%1,                Danielle%Ahn%NULL%1,                Mercedes%Martinez%NULL%1,                Kara Gross%Margolis%NULL%1]"
$ws.Range("E7").Value = "[Sonja A.%Rasmussen%skr9@ufl.edu%1,                John C.%Smulian%NULL%1,                John A.%Lednicky%NULL%1,                Tony S.%Wen%NULL%1,                Denise J.%Jamieson%NULL%1]"
$ws.Range("E8").Value = "[Elizabeth M.%Dufort%NULL%1,                Emilia H.%Koumans%NULL%1,                Eric J.%Chow%NULL%1,                Elizabeth M.%Rosenthal%NULL%2,                Elizabeth M.%Rosenthal%NULL%0,                Alison%Muse%NULL%2,                Alison%Muse%NULL%0,                Jemma%Rowlands%NULL%1,                Meredith A.%Barranco%NULL%1,                Angela M.%Maxted%NULL%1,                Eli S.%Rosenberg%NULL%1,                Delia%Easton%NULL%1,                Tomoko%Udo%NULL%1,                Jessica%Kumar%NULL%1,                Wendy%Pulver%NULL%1,                Lou%Smith%NULL%1,                Brad%Hutton%NULL%1,                Debra%Blog%NULL%1,                Howard%Zucker%NULL%1]"
$ws.Range("E9").Value = "[Shema%Hameed%shema.hameed@gstt.nhs.uk%1,                Heba%Elbaaly%NULL%2,                Heba%Elbaaly%NULL%0,                Catriona E. L.%Reid%NULL%2,                Catriona E. L.%Reid%NULL%0,                Rui M. F.%Santos%NULL%2,                Rui M. F.%Santos%NULL%0,                Vinay%Shivamurthy%NULL%2,                Vinay%Shivamurthy%NULL%0,                James%Wong%NULL%1,                K. Haran%Jogeesvaran%NULL%2,                K. Haran%Jogeesvaran%NULL%0]"
$ws.Range("E10").Value = "[Julie%Toubiana%NULL%1,                Clément%Poirault%NULL%1,                Alice%Corsia%NULL%1,                Fanny%Bajolle%NULL%2,                Jacques%Fourgeaud%NULL%1,                François%Angoulvant%NULL%1,                Agathe%Debray%NULL%1,                Romain%Basmaci%NULL%2,                Elodie%Salvador%NULL%1,                Sandra%Biscardi%NULL%1,                Pierre%Frange%NULL%1,                Martin%Chalumeau%NULL%1,                Jean-Laurent%Casanova%NULL%2,                Jean-Laurent%Casanova%NULL%0,                Jérémie F%Cohen%NULL%1,                Slimane%Allali%NULL%1]"
$ws.Range("E11").Value = "[Marie%Pouletty%NULL%1,                Charlotte%Borocco%NULL%1,                Naim%Ouldali%NULL%1,                Marion%Caseris%NULL%2,                Romain%Basmaci%NULL%0,                Noémie%Lachaume%NULL%1,                Philippe%Bensaid%NULL%1,                Samia%Pichard%NULL%1,                Hanane%Kouider%NULL%1,                Guillaume%Morelle%NULL%1,                Irina%Craiu%NULL%1,                Corinne%Pondarre%NULL%1,                Anna%Deho%NULL%1,                Arielle%Maroni%NULL%1,                Mehdi%Oualha%NULL%2,                Zahir%Amoura%NULL%1,                Julien%Haroche%NULL%1,                Juliette%Chommeloux%NULL%1,                Fanny%Bajolle%NULL%0,                Constance%Beyler%NULL%2,                Stéphane%Bonacorsi%NULL%1,                Guislaine%Carcelain%NULL%1,                Isabelle%Koné-Paut%NULL%1,                Brigitte%Bader-Meunier%NULL%2,                Brigitte%Bader-Meunier%NULL%0,                Albert%Faye%NULL%3,                Albert%Faye%NULL%0,                Ulrich%Meinzer%NULL%2,                Caroline%Galeotti%NULL%1,                Isabelle%Melki%NULL%2]"
$ws.Range("E12").Value = "[Zahra%Belhadjer%xref no email%1,        Mathilde%M\u00e9ot%xref no email%1,        Fanny%Bajolle%xref no email%1,        Diala%Khraiche%xref no email%1,        Antoine%Legendre%xref no email%1,        Samya%Abakka%xref no email%1,        Johanne%Auriau%xref no email%1,        Marion%Grimaud%xref no email%1,        Mehdi%Oualha%xref no email%1,        Maurice%Beghetti%xref no email%1,        Julie%Wacker%xref no email%1,        Caroline%Ovaert%xref no email%1,        Sebastien%Hascoet%xref no email%1,        Ma\u00eblle%Selegny%xref no email%1,        Sophie%Malekzadeh-Milani%xref no email%1,        Alice%Maltret%xref no email%1,        Gilles%Bosser%xref no email%1,        Nathan%Giroux%xref no email%1,        Laurent%Bonnemains%xref no email%1,        Jeanne%Bordet%xref no email%1,        Sylvie%Di Filippo%xref no email%1,        Pierre%Mauran%xref no email%1,        Sylvie%Falcon-Eicher%xref no email%1,        Jean-Beno\u00eet%Thambo%xref no email%1,        Bruno%Lefort%xref no email%1,        Pamela%Moceri%xref no email%1,        Lucile%Houyel%xref no email%1,        Sylvain%Renolleau%xref no email%1,        Damien%Bonnet%xref no email%1]"
$ws.Range("E13").Value = "[Tristan%Ramcharan%NULL%1,                Oscar%Nolan%NULL%1,                Chui Yi%Lai%NULL%1,                Nanda%Prabhu%NULL%1,                Raghu%Krishnamurthy%NULL%1,                Alex G.%Richter%NULL%1,                Deepthi%Jyothish%NULL%1,                Hari Krishnan%Kanthimathinathan%NULL%2,                Steven B.%Welch%NULL%1,                Scott%Hackett%NULL%1,                Eslam%Al-Abadi%NULL%1,                Barnaby R.%Scholefield%NULL%1,                Ashish%Chikermane%a.chikermane@nhs.net%1]"
$ws.Range("E14").Value = "[Marion%Grimaud%NULL%0,                Julie%Starck%NULL%1,                Michael%Levy%NULL%1,                Clémence%Marais%NULL%1,                Judith%Chareyre%NULL%1,                Diala%Khraiche%NULL%1,                Marianne%Leruez-Ville%NULL%1,                Pierre%Quartier%NULL%1,                Pierre Louis%Léger%NULL%1,                Guillaume%Geslain%NULL%1,                Nada%Semaan%NULL%1,                Florence%Moulin%NULL%1,                Matthieu%Bendavid%NULL%1,                Sandrine%Jean%NULL%1,                Géraldine%Poncelet%NULL%1,                Sylvain%Renolleau%NULL%1,                Mehdi%Oualha%mehdi.oualha@aphp.fr%0]"
$ws.Range("E15").Value = "[Mariawy%Riollano‐Cruz%NULL%1,                Esra%Akkoyun%NULL%1,                Eudys%Briceno‐Brito%NULL%1,                Shanna%Kowalsky%NULL%0,                James%Reed%NULL%1,                Roberto%Posada%NULL%1,                Emilia Mia%Sordillo%NULL%1,                Michael%Tosi%NULL%1,                Rebecca%Trachtman%NULL%0,                Alberto%Paniz‐Mondolfi%Alberto.Paniz-mondolfi@mountsinai.org%1]"
$ws.Range("E16").Value = "[Patrick%Davies%patrick.davies@nuh.nhs.uk%1,                Claire%Evans%NULL%1,                Hari Krishnan%Kanthimathinathan%NULL%0,                Jon%Lillie%NULL%1,                Joseph%Brierley%NULL%1,                Gareth%Waters%NULL%1,                Mae%Johnson%NULL%1,                Benedict%Griffiths%NULL%1,                Pascale%du Pré%NULL%1,                Zoha%Mohammad%NULL%1,                Akash%Deep%NULL%1,                Stephen%Playfor%NULL%1,                Davinder%Singh%NULL%1,                David%Inwald%NULL%1,                Michelle%Jardine%NULL%1,                Oliver%Ross%NULL%1,                Nayan%Shetty%NULL%1,                Mark%Worrall%NULL%1,                Ruchi%Sinha%NULL%1,                Ashwani%Koul%NULL%1,                Elizabeth%Whittaker%NULL%1,                Harish%Vyas%NULL%1,                Barnaby R%Scholefield%NULL%1,                Padmanabhan%Ramnarayan%NULL%1]"
$ws.Range("E17").Value = "[Naim%Ouldali%NULL%1,                Marie%Pouletty%NULL%1,                Patricia%Mariani%NULL%1,                Constance%Beyler%NULL%0,                Audrey%Blachier%NULL%1,                Stephane%Bonacorsi%NULL%1,                Kostas%Danis%NULL%1,                Maryline%Chomton%NULL%1,                Laure%Maurice%NULL%1,                Fleur%Le Bourgeois%NULL%1,                Marion%Caseris%NULL%0,                Jean%Gaschignard%NULL%1,                Julie%Poline%NULL%1,                Robert%Cohen%NULL%1,                Luigi%Titomanlio%NULL%1,                Albert%Faye%NULL%0,                Isabelle%Melki%NULL%0,                Ulrich%Meinzer%ulrich.meinzer@aphp.fr%0]"
$ws.Range("E18").Value = "[Leora R.%Feldstein%NULL%1,                Erica B.%Rose%NULL%1,                Steven M.%Horwitz%NULL%1,                Jennifer P.%Collins%NULL%1,                Margaret M.%Newhams%NULL%1,                Mary Beth F.%Son%NULL%1,                Jane W.%Newburger%NULL%1,                Lawrence C.%Kleinman%NULL%1,                Sabrina M.%Heidemann%NULL%1,                Amarilis A.%Martin%NULL%2,                Amarilis A.%Martin%NULL%0,                Aalok R.%Singh%NULL%1,                Simon%Li%NULL%1,                Keiko M.%Tarquinio%NULL%1,                Preeti%Jaggi%NULL%1,                Matthew E.%Oster%NULL%1,                Sheemon P.%Zackai%NULL%1,                Jennifer%Gillen%NULL%1,                Adam J.%Ratner%NULL%1,                Rowan F.%Walsh%NULL%1,                Julie C.%Fitzgerald%NULL%1,                Michael A.%Keenaghan%NULL%1,                Hussam%Alharash%NULL%1,                Sule%Doymaz%NULL%1,                Katharine N.%Clouser%NULL%1,                John S.%Giuliano%NULL%1,                Anjali%Gupta%NULL%1,                Robert M.%Parker%NULL%1,                Aline B.%Maddux%NULL%1,                Vinod%Havalad%NULL%1,                Stacy%Ramsingh%NULL%1,                Hulya%Bukulmez%NULL%1,                Tamara T.%Bradford%NULL%1,                Lincoln S.%Smith%NULL%1,                Mark W.%Tenforde%NULL%1,                Christopher L.%Carroll%NULL%1,                Becky J.%Riggs%NULL%1,                Shira J.%Gertz%NULL%1,                Ariel%Daube%NULL%1,                Amanda%Lansell%NULL%2,                Amanda%Lansell%NULL%0,                Alvaro%Coronado Munoz%NULL%1,                Charlotte V.%Hobbs%NULL%2,                Charlotte V.%Hobbs%NULL%0,                Kimberly L.%Marohn%NULL%1,                Natasha B.%Halasa%NULL%1,                Manish M.%Patel%NULL%1,                Adrienne G.%Randolph%NULL%1]"
